$wb = $excel.ActiveWorkbook

# Sheets: "Autos normal" (visible), "Aux" (hidden, contains list data)
$wsMain = $wb.Worksheets.Item("Autos normal")
$wsAux  = $wb.Worksheets.Item("Aux")

# --- Populate new helper columns on the Aux sheet ---
# Column A already has SI / NO (used by "aplica" defined name)
# Column B: A / D  -> used by new "A_D" defined name
# Column C: YES / NO -> used by new "YES_NO" defined name
$wsAux.Range("B1").Value = "A"
$wsAux.Range("B2").Value = "D"
$wsAux.Range("C1").Value = "YES"
$wsAux.Range("C2").Value = "NO"

# Set column B width to match the diff (stored OOXML width="11").
# The engine adds a fixed 5/6 character padding on top of ColumnWidth when
# serializing to the "width" attribute, so back it out here.
$wsAux.Columns.Item(2).ColumnWidth = 11 - (5/6)

# --- Add the new defined names ---
$wb.Names.Add("A_D", "=Aux!`$B`$1:`$B`$2")
$wb.Names.Add("YES_NO", "=Aux!`$C`$1:`$C`$2")

# --- Update selections / active cells to match the committed view state ---
$wsMain.Activate()
$wsMain.Range("E22").Select()

$wsAux.Activate()
$wsAux.Range("E11").Select()

# Return to the originally active sheet
$wsMain.Activate()
